# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Faerie Profits workbook sheets
# as described by the authoritative XML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 336494
$ws.Range("J17").Value = 336494
$ws.Range("L17").Value = 1009482
$ws.Range("N17").Value = -1009818
$ws.Range("H69").Value = 18642.857
$ws.Range("J69").Value = 18642.857
$ws.Range("L69").Value = 55928.571
$ws.Range("N69").Value = -57676.571
$ws.Range("H72").Value = 18642.857
$ws.Range("J72").Value = 18642.857
$ws.Range("L72").Value = 167785.713
$ws.Range("N72").Value = -176521.713
$ws.Range("H100").Value = 7571.72
$ws.Range("I100").Value = 940.5833
$ws.Range("K100").Value = 940.5833
$ws.Range("M100").Value = -399.5833
$ws.Range("H116").Value = 2914.2856
$ws.Range("I116").Value = 2900
$ws.Range("J116").Value = 2933.3333
$ws.Range("K116").Value = 2900
$ws.Range("L116").Value = 2933.3333
$ws.Range("M116").Value = 542
$ws.Range("N116").Value = -9817.3333
$ws.Range("H129").Value = 66667840
$ws.Range("I129").Value = 100000760
$ws.Range("K129").Value = 300002280
$ws.Range("M129").Value = -299997280
$ws.Range("H137").Value = 2154.074
$ws.Range("I137").Value = 2116.0435
$ws.Range("K137").Value = 6348.130500000001
$ws.Range("M137").Value = -3798.130500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 896.3333
$ws.Range("J4").Value = 890
$ws.Range("L4").Value = 890
$ws.Range("N4").Value = -1122
$ws.Range("H32").Value = 6275.125
$ws.Range("I32").Value = 5723.125
$ws.Range("K32").Value = 5723.125
$ws.Range("M32").Value = -5436.125
$ws.Range("H61").Value = 3131.2285
$ws.Range("I61").Value = 2559.3333
$ws.Range("J61").Value = 6562.6
$ws.Range("K61").Value = 2559.3333
$ws.Range("L61").Value = 6562.6
$ws.Range("M61").Value = -2347.3333
$ws.Range("N61").Value = -6986.6
$ws.Range("H74").Value = 2465.4443
$ws.Range("I74").Value = 1773.75
$ws.Range("K74").Value = 1773.75
$ws.Range("M74").Value = -899.75
$ws.Range("H77").Value = 2465.4443
$ws.Range("I77").Value = 1773.75
$ws.Range("K77").Value = 8868.75
$ws.Range("M77").Value = -4500.75
$ws.Range("H97").Value = 2123.5789
$ws.Range("I97").Value = 1424.625
$ws.Range("K97").Value = 1424.625
$ws.Range("M97").Value = -928.625
$ws.Range("H101").Value = 75000
$ws.Range("J101").Value = 75000
$ws.Range("L101").Value = 75000
$ws.Range("N101").Value = -81490
$ws.Range("H102").Value = 3213.8333
$ws.Range("I102").Value = 2316.9
$ws.Range("K102").Value = 2316.9
$ws.Range("M102").Value = -694.9000000000001
$ws.Range("H110").Value = 2345.0833
$ws.Range("I110").Value = 1268.625
$ws.Range("K110").Value = 1268.625
$ws.Range("M110").Value = 776.375
$ws.Range("H122").Value = 1458
$ws.Range("I122").Value = 1342.4546
$ws.Range("K122").Value = 4027.3638
$ws.Range("M122").Value = -1577.3638
$ws.Range("H136").Value = 3131.2285
$ws.Range("I136").Value = 2559.3333
$ws.Range("J136").Value = 6562.6
$ws.Range("K136").Value = 7677.999899999999
$ws.Range("L136").Value = 19687.8
$ws.Range("M136").Value = -5127.999899999999
$ws.Range("N136").Value = -24787.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3836.2
$ws.Range("I94").Value = 3199.5
$ws.Range("J94").Value = 4791.25
$ws.Range("K94").Value = 3199.5
$ws.Range("L94").Value = 4791.25
$ws.Range("M94").Value = -2748.5
$ws.Range("N94").Value = -5693.25
$ws.Range("H105").Value = 4414.143
$ws.Range("I105").Value = 3167.875
$ws.Range("J105").Value = 8402.200000000001
$ws.Range("K105").Value = 3167.875
$ws.Range("L105").Value = 8402.200000000001
$ws.Range("M105").Value = -1420.875
$ws.Range("N105").Value = -11896.2
$ws.Range("H107").Value = 10870650
$ws.Range("J107").Value = 1703
$ws.Range("L107").Value = 1703
$ws.Range("N107").Value = -5543

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 568.1
$ws.Range("I22").Value = 543.3077
$ws.Range("J22").Value = 614.1429000000001
$ws.Range("K22").Value = 543.3077
$ws.Range("L22").Value = 614.1429000000001
$ws.Range("M22").Value = -193.3077
$ws.Range("N22").Value = -1314.1429
$ws.Range("H31").Value = 3165.3462
$ws.Range("I31").Value = 1566.0834
$ws.Range("J31").Value = 4536.143
$ws.Range("K31").Value = 1566.0834
$ws.Range("L31").Value = 4536.143
$ws.Range("M31").Value = -1271.0834
$ws.Range("N31").Value = -5126.143
$ws.Range("H34").Value = 3165.3462
$ws.Range("I34").Value = 1566.0834
$ws.Range("J34").Value = 4536.143
$ws.Range("K34").Value = 1566.0834
$ws.Range("L34").Value = 4536.143
$ws.Range("M34").Value = -1364.0834
$ws.Range("N34").Value = -4940.143
$ws.Range("H43").Value = 61315.332
$ws.Range("J43").Value = 61315.332
$ws.Range("L43").Value = 61315.332
$ws.Range("N43").Value = -61683.332
$ws.Range("H58").Value = 2134.4285
$ws.Range("I58").Value = 1918.2667
$ws.Range("J58").Value = 2674.8333
$ws.Range("K58").Value = 1918.2667
$ws.Range("L58").Value = 2674.8333
$ws.Range("M58").Value = -1715.2667
$ws.Range("N58").Value = -3080.8333
$ws.Range("H96").Value = 99999
$ws.Range("J96").Value = 99999
$ws.Range("L96").Value = 99999
$ws.Range("N96").Value = -105491
$ws.Range("H101").Value = 61315.332
$ws.Range("J101").Value = 61315.332
$ws.Range("L101").Value = 61315.332
$ws.Range("N101").Value = -67805.33199999999
$ws.Range("H122").Value = 3644.238
$ws.Range("I122").Value = 3460.111
$ws.Range("K122").Value = 10380.333
$ws.Range("M122").Value = -7930.332999999999
$ws.Range("H136").Value = 2134.4285
$ws.Range("I136").Value = 1918.2667
$ws.Range("J136").Value = 2674.8333
$ws.Range("K136").Value = 5754.800099999999
$ws.Range("L136").Value = 8024.499899999999
$ws.Range("M136").Value = -3204.800099999999
$ws.Range("N136").Value = -13124.4999
$ws.Range("H141").Value = 121492.445
$ws.Range("J141").Value = 131679
$ws.Range("L141").Value = 131679
$ws.Range("N141").Value = -142039

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 636.7143
$ws.Range("J122").Value = 973.8570999999999
$ws.Range("L122").Value = 8764.713899999999
$ws.Range("N122").Value = -13664.7139
$ws.Range("H131").Value = 1789.2325
$ws.Range("I131").Value = 1182
$ws.Range("J131").Value = 1834.775
$ws.Range("K131").Value = 3546
$ws.Range("L131").Value = 5504.325000000001
$ws.Range("M131").Value = 1494
$ws.Range("N131").Value = -15584.325
$ws.Range("H132").Value = 3365.875
$ws.Range("I132").Value = 1049.75
$ws.Range("K132").Value = 9447.75
$ws.Range("M132").Value = -6917.75
$ws.Range("H138").Value = 4786567.5
$ws.Range("I138").Value = 670670.25
$ws.Range("J138").Value = 12503875
$ws.Range("K138").Value = 2012010.75
$ws.Range("L138").Value = 37511625
$ws.Range("M138").Value = -2006870.75
$ws.Range("N138").Value = -37521905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4740.7144
$ws.Range("I70").Value = 4796
$ws.Range("J70").Value = 4685.4287
$ws.Range("K70").Value = 4796
$ws.Range("L70").Value = 4685.4287
$ws.Range("M70").Value = -4526
$ws.Range("N70").Value = -5225.4287
$ws.Range("H73").Value = 4740.7144
$ws.Range("I73").Value = 4796
$ws.Range("J73").Value = 4685.4287
$ws.Range("K73").Value = 4796
$ws.Range("L73").Value = 4685.4287
$ws.Range("M73").Value = -3860
$ws.Range("N73").Value = -6557.4287
$ws.Range("H132").Value = 6947147
$ws.Range("I132").Value = 8336016.5
$ws.Range("J132").Value = 2801.25
$ws.Range("K132").Value = 25008049.5
$ws.Range("L132").Value = 8403.75
$ws.Range("M132").Value = -25005519.5
$ws.Range("N132").Value = -13463.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2654.85
$ws.Range("I40").Value = 2746.8823
$ws.Range("J40").Value = 2133.3333
$ws.Range("K40").Value = 2746.8823
$ws.Range("L40").Value = 2133.3333
$ws.Range("M40").Value = -2610.8823
$ws.Range("N40").Value = -2405.3333
$ws.Range("H46").Value = 6288.8887
$ws.Range("I46").Value = 2673.4546
$ws.Range("J46").Value = 8774.5
$ws.Range("K46").Value = 2673.4546
$ws.Range("L46").Value = 8774.5
$ws.Range("M46").Value = -2485.4546
$ws.Range("N46").Value = -9150.5
$ws.Range("H100").Value = 3947.3
$ws.Range("I100").Value = 3448
$ws.Range("K100").Value = 3448
$ws.Range("M100").Value = -2907
$ws.Range("H111").Value = 79900
$ws.Range("J111").Value = 79900
$ws.Range("L111").Value = 79900
$ws.Range("N111").Value = -88080
$ws.Range("H132").Value = 2775.9614
$ws.Range("I132").Value = 2782.524
$ws.Range("J132").Value = 2748.4
$ws.Range("K132").Value = 8347.572
$ws.Range("L132").Value = 8245.200000000001
$ws.Range("M132").Value = -5817.572
$ws.Range("N132").Value = -13305.2
$ws.Range("H140").Value = 133690.5
$ws.Range("J140").Value = 133690.5
$ws.Range("L140").Value = 133690.5
$ws.Range("N140").Value = -144050.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1916.1786
$ws.Range("J100").Value = 3579.6
$ws.Range("L100").Value = 7159.2
$ws.Range("N100").Value = -8241.200000000001
$ws.Range("H107").Value = 522.625
$ws.Range("I107").Value = 449.5
$ws.Range("J107").Value = 644.5
$ws.Range("K107").Value = 1348.5
$ws.Range("L107").Value = 1933.5
$ws.Range("M107").Value = 571.5
$ws.Range("N107").Value = -5773.5
$ws.Range("H113").Value = 5953304
$ws.Range("I113").Value = 11905368
$ws.Range("J113").Value = 1240.2858
$ws.Range("K113").Value = 35716104
$ws.Range("L113").Value = 3720.8574
$ws.Range("M113").Value = -35713934
$ws.Range("N113").Value = -8060.857400000001
$ws.Range("H122").Value = 2920
$ws.Range("I122").Value = 2807.4348
$ws.Range("J122").Value = 3955.6
$ws.Range("K122").Value = 8422.304400000001
$ws.Range("L122").Value = 11866.8
$ws.Range("M122").Value = -5972.304400000001
$ws.Range("N122").Value = -16766.8

